# Re-worked to use modular error message creation and includes manifests for testing
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cell E8 now holds the new shared-string value "123dgfg" instead of the
# numeric literal 123456.
$ws.Range("E8").Value = "123dgfg"

# The saved selection moves from B2:B8 to the single cell E8.
$ws.Range("E8").Select()

# Reflect the updated window geometry from the saved workbook view.
$win = $excel.ActiveWindow
$win.Left = 4100
$win.Top = -20880
$win.Width = 13900
$win.Height = 16500
